$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 8
$ws.Range("I11").Value = 8
$ws.Range("K11").Value = 8
$ws.Range("M11").Value = 132

$ws.Range("H40").Value = 1977946.5
$ws.Range("I40").Value = 1070
$ws.Range("J40").Value = 3107590.2
$ws.Range("K40").Value = 1070
$ws.Range("L40").Value = 3107590.2
$ws.Range("M40").Value = -895
$ws.Range("N40").Value = -3107940.2

$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H64").Value = 4800
$ws.Range("J64").Value = 5000
$ws.Range("L64").Value = 5000
$ws.Range("N64").Value = -5496

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H67").Value = 4800
$ws.Range("J67").Value = 5000
$ws.Range("L67").Value = 5000
$ws.Range("N67").Value = -6716

$ws.Range("H127").Value = 1266
$ws.Range("J127").Value = 1501.4445
$ws.Range("L127").Value = 4504.333500000001
$ws.Range("N127").Value = -14424.3335

$ws.Range("H135").Value = 20839874
$ws.Range("I135").Value = 628.4737
$ws.Range("K135").Value = 5656.263300000001
$ws.Range("M135").Value = -3121.263300000001

$ws.Range("H137").Value = 21115.156
$ws.Range("I137").Value = 1284.475
$ws.Range("J137").Value = 93226.73
$ws.Range("K137").Value = 3853.425
$ws.Range("L137").Value = 279680.19
$ws.Range("M137").Value = -1303.425
$ws.Range("N137").Value = -284780.19

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 28685.205
$ws.Range("I32").Value = 32824.656
$ws.Range("K32").Value = 32824.656
$ws.Range("M32").Value = -32537.656

$ws.Range("H45").Value = 2721.3125
$ws.Range("I45").Value = 2454.72
$ws.Range("K45").Value = 2454.72
$ws.Range("M45").Value = -2077.72

$ws.Range("H61").Value = 1731.6923
$ws.Range("I61").Value = 1383.6522
$ws.Range("J61").Value = 4400
$ws.Range("K61").Value = 1383.6522
$ws.Range("L61").Value = 4400
$ws.Range("M61").Value = -1171.6522
$ws.Range("N61").Value = -4824

$ws.Range("H63").Value = 2228.5
$ws.Range("I63").Value = 2251
$ws.Range("K63").Value = 2251
$ws.Range("M63").Value = -1565

$ws.Range("H66").Value = 2228.5
$ws.Range("I66").Value = 2251
$ws.Range("K66").Value = 11255
$ws.Range("M66").Value = -7823

$ws.Range("H74").Value = 76924060
$ws.Range("I74").Value = 76924060
$ws.Range("K74").Value = 76924060
$ws.Range("M74").Value = -76923186

$ws.Range("H77").Value = 76924060
$ws.Range("I77").Value = 76924060
$ws.Range("K77").Value = 384620300
$ws.Range("M77").Value = -384615932

$ws.Range("H136").Value = 1731.6923
$ws.Range("I136").Value = 1383.6522
$ws.Range("J136").Value = 4400
$ws.Range("K136").Value = 4150.9566
$ws.Range("L136").Value = 13200
$ws.Range("M136").Value = -1600.9566
$ws.Range("N136").Value = -18300

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 3833
$ws.Range("I11").Value = 5249.5
$ws.Range("K11").Value = 5249.5
$ws.Range("M11").Value = -5109.5

$ws.Range("H22").Value = 1000473.9
$ws.Range("I22").Value = 1428920.1
$ws.Range("J22").Value = 766
$ws.Range("K22").Value = 1428920.1
$ws.Range("L22").Value = 766
$ws.Range("M22").Value = -1428747.1
$ws.Range("N22").Value = -1112

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15837.04
$ws.Range("I31").Value = 19913.295
$ws.Range("J31").Value = 7175
$ws.Range("K31").Value = 19913.295
$ws.Range("L31").Value = 7175
$ws.Range("M31").Value = -19618.295
$ws.Range("N31").Value = -7765

$ws.Range("H34").Value = 15837.04
$ws.Range("I34").Value = 19913.295
$ws.Range("J34").Value = 7175
$ws.Range("K34").Value = 19913.295
$ws.Range("L34").Value = 7175
$ws.Range("M34").Value = -19711.295
$ws.Range("N34").Value = -7579

$ws.Range("H99").Value = 23814068
$ws.Range("I99").Value = 4260
$ws.Range("J99").Value = 83338584
$ws.Range("K99").Value = 4260
$ws.Range("L99").Value = 83338584
$ws.Range("M99").Value = -2762
$ws.Range("N99").Value = -83341580

$ws.Range("H126").Value = 23814068
$ws.Range("I126").Value = 4260
$ws.Range("J126").Value = 83338584
$ws.Range("K126").Value = 12780
$ws.Range("L126").Value = 250015752
$ws.Range("M126").Value = -10310
$ws.Range("N126").Value = -250020692

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 3219.6667
$ws.Range("I58").Value = 1000
$ws.Range("J58").Value = 3663.6
$ws.Range("K58").Value = 3000
$ws.Range("L58").Value = 10990.8
$ws.Range("M58").Value = -2872
$ws.Range("N58").Value = -11246.8

$ws.Range("H122").Value = 1125.8235
$ws.Range("J122").Value = 1370.6923
$ws.Range("L122").Value = 12336.2307
$ws.Range("N122").Value = -17236.2307

$ws.Range("H131").Value = 764.6
$ws.Range("I131").Value = 330
$ws.Range("J131").Value = 818.31464
$ws.Range("K131").Value = 990
$ws.Range("L131").Value = 2454.94392
$ws.Range("M131").Value = 4050
$ws.Range("N131").Value = -12534.94392

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 29990
$ws.Range("J39").Value = 29990
$ws.Range("L39").Value = 29990
$ws.Range("N39").Value = -31054

$ws.Range("H102").Value = 21741472
$ws.Range("I102").Value = 29414408
$ws.Range("J102").Value = 1485.6666
$ws.Range("K102").Value = 29414408
$ws.Range("L102").Value = 1485.6666
$ws.Range("M102").Value = -29412786
$ws.Range("N102").Value = -4729.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2321
$ws.Range("I68").Value = 2323.2307
$ws.Range("J68").Value = 2316.1667
$ws.Range("K68").Value = 2323.2307
$ws.Range("L68").Value = 2316.1667
$ws.Range("M68").Value = -1574.2307
$ws.Range("N68").Value = -3814.1667

$ws.Range("H71").Value = 2321
$ws.Range("I71").Value = 2323.2307
$ws.Range("J71").Value = 2316.1667
$ws.Range("K71").Value = 11616.1535
$ws.Range("L71").Value = 11580.8335
$ws.Range("M71").Value = -7872.1535
$ws.Range("N71").Value = -19068.8335

$ws.Range("H132").Value = 1195.625
$ws.Range("I132").Value = 884.0833
$ws.Range("K132").Value = 2652.2499
$ws.Range("M132").Value = -122.2498999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3497138.5
$ws.Range("I107").Value = 873.5714
$ws.Range("J107").Value = 7576114.5
$ws.Range("K107").Value = 2620.7142
$ws.Range("L107").Value = 22728343.5
$ws.Range("M107").Value = -700.7142000000003
$ws.Range("N107").Value = -22732183.5

$ws.Range("H122").Value = 1272.1875
$ws.Range("I122").Value = 1114.5454
$ws.Range("J122").Value = 1619
$ws.Range("K122").Value = 3343.6362
$ws.Range("L122").Value = 4857
$ws.Range("M122").Value = -893.6361999999999
$ws.Range("N122").Value = -9757

$ws.Range("H132").Value = 1100.2307
$ws.Range("I132").Value = 765.7826
$ws.Range("J132").Value = 3664.3333
$ws.Range("K132").Value = 2297.3478
$ws.Range("L132").Value = 10992.9999
$ws.Range("M132").Value = 232.6522
$ws.Range("N132").Value = -16052.9999

$ws.Range("H136").Value = 27779200
$ws.Range("I136").Value = 30304264
$ws.Range("K136").Value = 90912792
$ws.Range("M136").Value = -90910242
